$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.11095929145813
$ws.Range("B1").Value = 2.576983213424683
$ws.Range("C1").Value = 2.694610595703125
$ws.Range("D1").Value = 3.195507764816284
$ws.Range("E1").Value = 0.7914511561393738
